# "finish dev of sell items"
#
# PlayerInitData.xlsx changes:
#   1. intInit: add two new int-constant rows used to configure the new
#      "sell items for a cart" feature - resourcesPerCart=1000,
#      materialsPerCart=1.
#   2. resources: drop the now-unused INT_energy column (energy was removed
#      from the sellable/collectible resource list), shifting INT_blood
#      left into its place.
#   3. Minor selection/view bookkeeping to mirror the authored workbook.

$wb = $excel.ActiveWorkbook

$intInit   = $wb.Worksheets.Item("intInit")
$resources = $wb.Worksheets.Item("resources")

# --- intInit: new rows for the cart-selling feature -----------------------
$intInit.Range("A3").Value = "resourcesPerCart"
$intInit.Range("B3").Value = 1000
$intInit.Range("A4").Value = "materialsPerCart"
$intInit.Range("B4").Value = 1

# --- resources: remove the INT_energy column (J), INT_blood shifts left ---
$resources.Range("J1:J4").Delete(-4159)   # xlShiftToLeft

# --- view/selection bookkeeping --------------------------------------------
$resources.Range("J5").Select()
$intInit.Range("B5").Select()
$intInit.Activate()
